$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 ("btn quit"), shifting every row below it down by one.
$ws.Rows.Item(4).Insert()

# Copy the formatting (style + row height) of the row now below the inserted one
# so the new row matches the look of the other data rows.
$ws.Range("A5:F5").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)

# Populate the new row with the "Log out" button localization entry.
$ws.Range("A4").Value = $ws.Range("A5").Value()
$ws.Range("B4").Value = "btn logout"
$ws.Range("C4").Value = "Log out"
$ws.Range("D4").Value = "Déconnexion"

# Update the active selection to match the saved view state.
$ws.Range("B5").Select()

Write-Output "applied logout row insertion"
